# Update "想去人数" (interest count, column F) values across all four sheets
# to reflect the refreshed data snapshot (gh-pages output generated at 456a3b4).
# Every other cell/column is left untouched.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1748
$ws.Range("F5").Value = 457
$ws.Range("F7").Value = 247
$ws.Range("F8").Value = 1198
$ws.Range("F11").Value = 881
$ws.Range("F12").Value = 688
$ws.Range("F14").Value = 509
$ws.Range("F15").Value = 141
$ws.Range("F17").Value = 174
$ws.Range("F18").Value = 2925
$ws.Range("F19").Value = 2624
$ws.Range("F26").Value = 5278
$ws.Range("F28").Value = 986
$ws.Range("F31").Value = 312
$ws.Range("F32").Value = 1098
$ws.Range("F35").Value = 288

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1122
$ws.Range("F19").Value = 43
$ws.Range("F24").Value = 317
$ws.Range("F26").Value = 3934
$ws.Range("F30").Value = 199
$ws.Range("F33").Value = 166

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2459
$ws.Range("F9").Value = 1324
$ws.Range("F10").Value = 360

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2459
$ws.Range("F5").Value = 1748
$ws.Range("F7").Value = 1324
$ws.Range("F8").Value = 360
$ws.Range("F11").Value = 457
$ws.Range("F13").Value = 247
$ws.Range("F14").Value = 1198
$ws.Range("F16").Value = 881
$ws.Range("F17").Value = 688
$ws.Range("F18").Value = 1122
$ws.Range("F19").Value = 1122
$ws.Range("F21").Value = 509
$ws.Range("F23").Value = 174
$ws.Range("F24").Value = 2925
$ws.Range("F25").Value = 2625
$ws.Range("F30").Value = 5278
$ws.Range("F32").Value = 986
$ws.Range("F38").Value = 312
$ws.Range("F41").Value = 43
$ws.Range("F44").Value = 317
$ws.Range("F45").Value = 317
$ws.Range("F46").Value = 1098
$ws.Range("F47").Value = 199
$ws.Range("F49").Value = 166
$ws.Range("F51").Value = 288

$wb.Save()
